$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.710.60'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '1.722.14'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('D4').Value = '0.9981'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '239.78'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').Value = '0.4822'
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('D8').Value = '0.2575'
$ws.Range('E8').Value = '  -0.41%  '
$ws.Range('D9').Value = '0.06182'
$ws.Range('D10').Value = '1.724.20'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('E11').Value = '  +2.66%  '
$ws.Range('D12').Value = '0.06849'
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('D13').Value = '0.6038'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').Value = '4.459'
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').Value = '76.78'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '26.538.47'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').Value = '0.9982'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').Value = '0.000007142'
$ws.Range('E19').Value = '  -1.41%  '
$ws.Range('D20').Value = '11.34'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').Value = '1.946.14'
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('E22').Value = '  -1.01%  '
$ws.Range('D23').Value = '8.567'
$ws.Range('D24').Value = '5.052'
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('D25').Value = '139.24'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('D26').Value = '15.20'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('D27').Value = '1.767'
$ws.Range('E27').Value = '  +2.52%  '
$ws.Range('D28').Value = '105.95'
$ws.Range('E28').Value = '  -0.82%  '
$ws.Range('E29').Value = '  -2.59%  '
$ws.Range('D30').Value = '4.025'
$ws.Range('E30').Value = '  +1.90%  '
$ws.Range('D31').Value = '0.07901'
$ws.Range('E31').Value = '  -0.89%  '
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('D33').Value = '0.04476'
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('D34').Value = '0.9979'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D36').Value = '0.9981'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '0.6161'
$ws.Range('E37').Value = '  -1.69%  '
$ws.Range('D38').Value = '0.9250'
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('E39').Value = '  +2.76%  '
$ws.Range('D40').Value = '2.437'
$ws.Range('E40').Value = '  +2.06%  '
$ws.Range('D41').Value = '0.9980'
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').Value = '0.01488'
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('D43').Value = '5.598'
$ws.Range('E43').Value = '  +4.93%  '
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('D45').Value = '0.3821'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').Value = '6.771'
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('E47').Value = '  -1.46%  '
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('D49').Value = '7.837'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('D50').Value = '30.00'
$ws.Range('D51').Value = '1.234'
$ws.Range('E51').Value = '  +0.56%  '
